$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price ("D") cells hold numeric-looking text (e.g. "214.12", "26.012.39").
# A plain .Value assignment lets Excel auto-coerce those into real numbers,
# so force the cell to Text format before writing, then clear the formatting
# again (the source cells carry no explicit style) to avoid leaving a stray
# number-format style behind while keeping the stored value as text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.012.39"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.633.93"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  -1.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0623"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.52"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.860.94"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.22"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.641.89"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0747"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.010.26"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.93"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "190.71"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("E22").Value = "  -3.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.15"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.41"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.75"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("E29").Value = "  -1.85%  "
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0485"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.84%  "
$ws.Range("E32").Value = "  -2.62%  "
$ws.Range("E33").Value = "  -3.88%  "
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.875"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.131.48"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.84"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.771.12"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0113"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.34"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.24%  "
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.48"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("E51").Value = "  +0.13%  "
